$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lower-case the header row text
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "subject1"
$ws.Range("C1").Value = "subject2"
$ws.Range("D1").Value = "subject3"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "grade"

# Append new rows of student data
$data = @(
    @("sai", 9, 90, 99, 198, "C"),
    @("Rani", 99, 99, 99, 297, "A"),
    @("rahul", 11, 11, 11, 33, "Fail"),
    @("Rahul", 11, 12, 13, 36, "Fail"),
    @("rahul", 1, 2, 3, 6, "Fail"),
    @("Sai", 100, 100, 100, 300, "A")
)

$row = 6
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $row++
}
